$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 840, shifting rows 840:883 down to 841:884.
$ws.Rows.Item(840).Insert()

# Populate the newly inserted row 840 with the new weekly observation
# (same Mercado/Region/Variedad/Calidad as the former row 840, new Fecha + Volumen).
$ws.Range("A840").Value = 10
$ws.Range("B840").Value = "Vega Modelo de Temuco"
$ws.Range("C840").Value = "La Araucanía"
$ws.Range("D840").Value = 45147
$ws.Range("E840").Value = 9
$ws.Range("F840").Value = 100112045
$ws.Range("G840").Value = "Zapallo"
$ws.Range("H840").Value = "Camote"
$ws.Range("I840").Value = "1a (guarda)"
$ws.Range("J840").Value = 125
$ws.Range("K840").Value = 600
$ws.Range("L840").Value = 600
$ws.Range("M840").Value = 600
$ws.Range("N840").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O840").Value = "Región del Maule"
$ws.Range("P840").Value = 600
$ws.Range("Q840").Value = 1
$ws.Range("R840").Value = "Hortaliza"

# Match the date-number style used by the rest of column D (numFmt 165).
$ws.Range("D840").NumberFormat = $ws.Range("D841").NumberFormat
